# Update the email address on the "Registration" sheet (cell A2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")
$ws.Range("A2").Value = "zqio@test.com"
